$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scouting Admin")

# Update the "Manage Season" section sync/date stamps (rows 6-14) from
# 45317 -> 45412 (B column)
$ws.Range("B6:B14").Value = 45412

# Row 14 "Add Event" becomes "Add TBA Code Event" (text edit in place)
$ws.Range("A14").Value = "Add TBA Code Event"

# Insert a brand-new row right after row 14 for the new
# "Add Event manual input" feature entry, pushing everything below down
# by one (mirrors the row-14-copy-down formatting Excel performs on a
# manual Insert).
$ws.Rows(15).Insert()
$ws.Range("A15").Value = "Add Event manual input"
$ws.Range("B15").Value = 45412

# The old "Add Team" / "Link Team To Event" / "Remove Team From Event"
# rows (now at 16/17/18 after the insert above) lose their date stamp
# and gain a "Need to fix" note in column C.
$ws.Range("B16:B18").ClearContents()
$ws.Range("C16").Value = "Need to fix"
$ws.Range("C17").Value = "Need to fix"
$ws.Range("C18").Value = "Need to fix"

# Refresh the sheet's recorded selection to match the author's final
# cursor position.
$ws.Range("D18").Select() | Out-Null
